$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# G2: append "/3" to the existing formula (divide CDM sum by 3)
$ws.Range("G2").Formula = "=SUM(ROUND(ABS(C2/SQRT(B2^2+C2^2-2*C2*B2*ROUND(COS(IF(D2<90,RADIANS(180-D2),RADIANS(D2))),4))-1/SQRT(2)),4),ROUND(ABS(B2/SQRT(C2^2+B2^2-2*B2*C2*ROUND(COS(IF(D2<90,RADIANS(180-D2),RADIANS(D2))),4))-1/SQRT(2)),4),ROUND(ABS(C2/SQRT(A2^2+C2^2-2*C2*A2*ROUND(COS(IF(E2<90,RADIANS(180-E2),RADIANS(E2))),4))-1/SQRT(2)),4),ROUND(ABS(A2/SQRT(C2^2+A2^2-2*A2*C2*ROUND(COS(IF(E2<90,RADIANS(180-E2),RADIANS(E2))),4))-1/SQRT(2)),4),ROUND(ABS(A2/SQRT(A2^2+B2^2-2*A2*B2*ROUND(COS(IF(F2<90,RADIANS(180-F2),RADIANS(F2))),4))-1/SQRT(2)),4),ROUND(ABS(B2/SQRT(A2^2+B2^2-2*A2*B2*ROUND(COS(IF(F2<90,RADIANS(180-F2),RADIANS(F2))),4))-1/SQRT(2)),4))/3"

# G3:G11 shared formula group: extend to the full 6-term CDM expression, divided by 3
$ws.Range("G3").Formula = "=SUM(ROUND(ABS(C3/SQRT(B3^2+C3^2-2*C3*B3*ROUND(COS(IF(D3<90,RADIANS(180-D3),RADIANS(D3))),4))-1/SQRT(2)),4),ROUND(ABS(B3/SQRT(C3^2+B3^2-2*B3*C3*ROUND(COS(IF(D3<90,RADIANS(180-D3),RADIANS(D3))),4))-1/SQRT(2)),4),ROUND(ABS(C3/SQRT(A3^2+C3^2-2*C3*A3*ROUND(COS(IF(E3<90,RADIANS(180-E3),RADIANS(E3))),4))-1/SQRT(2)),4),ROUND(ABS(A3/SQRT(C3^2+A3^2-2*A3*C3*ROUND(COS(IF(E3<90,RADIANS(180-E3),RADIANS(E3))),4))-1/SQRT(2)),4),ROUND(ABS(A3/SQRT(A3^2+B3^2-2*A3*B3*ROUND(COS(IF(F3<90,RADIANS(180-F3),RADIANS(F3))),4))-1/SQRT(2)),4),ROUND(ABS(B3/SQRT(A3^2+B3^2-2*A3*B3*ROUND(COS(IF(F3<90,RADIANS(180-F3),RADIANS(F3))),4))-1/SQRT(2)),4))/3"
$ws.Range("G4:G11").Formula = "=SUM(ROUND(ABS(C4/SQRT(B4^2+C4^2-2*C4*B4*ROUND(COS(IF(D4<90,RADIANS(180-D4),RADIANS(D4))),4))-1/SQRT(2)),4),ROUND(ABS(B4/SQRT(C4^2+B4^2-2*B4*C4*ROUND(COS(IF(D4<90,RADIANS(180-D4),RADIANS(D4))),4))-1/SQRT(2)),4),ROUND(ABS(C4/SQRT(A4^2+C4^2-2*C4*A4*ROUND(COS(IF(E4<90,RADIANS(180-E4),RADIANS(E4))),4))-1/SQRT(2)),4),ROUND(ABS(A4/SQRT(C4^2+A4^2-2*A4*C4*ROUND(COS(IF(E4<90,RADIANS(180-E4),RADIANS(E4))),4))-1/SQRT(2)),4),ROUND(ABS(A4/SQRT(A4^2+B4^2-2*A4*B4*ROUND(COS(IF(F4<90,RADIANS(180-F4),RADIANS(F4))),4))-1/SQRT(2)),4),ROUND(ABS(B4/SQRT(A4^2+B4^2-2*A4*B4*ROUND(COS(IF(F4<90,RADIANS(180-F4),RADIANS(F4))),4))-1/SQRT(2)),4))/3"

# Update the active selection shown in the sheet view
$ws.Range("P22").Select()
